$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1 = 14, Q1 = 15, matching the existing bold/border style ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").Borders.LineStyle = 1
$ws.Range("P1:Q1").HorizontalAlignment = -4108
$ws.Range("P1:Q1").VerticalAlignment = -4160

# --- Data rows 2-25: swap I<->(new value), K<->M, O<->(new) per the diff, and add P/Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column = 2
}
